$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data from the Feb 16 2023 GitHub Actions refresh.
# Each entry forces the literal text (leading apostrophe) so the numeric-looking
# strings are not reinterpreted as Excel numbers/percentages, matching the sheet's
# existing text-valued Price/Volume columns.
$updates = @(
    @{ Cell = "D2"; Value = "320.59" }
    @{ Cell = "E2"; Value = "7.70%" }
    @{ Cell = "D3"; Value = "48.06" }
    @{ Cell = "E3"; Value = "15.23%" }
    @{ Cell = "D4"; Value = "5.256" }
    @{ Cell = "E4"; Value = "4.90%" }
    @{ Cell = "D5"; Value = "0.08113" }
    @{ Cell = "E5"; Value = "7.85%" }
    @{ Cell = "D6"; Value = "4.606" }
    @{ Cell = "E6"; Value = "5.30%" }
    @{ Cell = "D7"; Value = "1.651" }
    @{ Cell = "E7"; Value = "3.03%" }
    @{ Cell = "D8"; Value = "1.189" }
    @{ Cell = "E8"; Value = "29.70%" }
    @{ Cell = "E9"; Value = "11.74%" }
    @{ Cell = "D10"; Value = "0.1944" }
    @{ Cell = "E10"; Value = "6.23%" }
    @{ Cell = "D11"; Value = "0.09578" }
    @{ Cell = "E11"; Value = "6.82%" }
    @{ Cell = "D12"; Value = "0.04599" }
    @{ Cell = "E12"; Value = "12.23%" }
    @{ Cell = "E13"; Value = "0.02%" }
    @{ Cell = "D14"; Value = "0.001338" }
    @{ Cell = "E14"; Value = "4.78%" }
    @{ Cell = "D15"; Value = "0.005938" }
    @{ Cell = "E15"; Value = "-0.75%" }
    @{ Cell = "D16"; Value = "3.364" }
    @{ Cell = "E16"; Value = "0.76%" }
    @{ Cell = "D17"; Value = "2.436" }
    @{ Cell = "E17"; Value = "1.43%" }
    @{ Cell = "D18"; Value = "0.3396" }
    @{ Cell = "E18"; Value = "2.02%" }
    @{ Cell = "D19"; Value = "8.161" }
    @{ Cell = "E19"; Value = "-1.55%" }
    @{ Cell = "D20"; Value = "0.1410" }
    @{ Cell = "E20"; Value = "4.34%" }
    @{ Cell = "D21"; Value = "0.3149" }
    @{ Cell = "E21"; Value = "1.49%" }
    @{ Cell = "D22"; Value = "0.04288" }
    @{ Cell = "E22"; Value = "4.54%" }
    @{ Cell = "D23"; Value = "0.001306" }
    @{ Cell = "E23"; Value = "3.14%" }
    @{ Cell = "D24"; Value = "0.004248" }
    @{ Cell = "E24"; Value = "9.11%" }
    @{ Cell = "E25"; Value = "3.77%" }
    @{ Cell = "D26"; Value = "0.0003539" }
    @{ Cell = "E26"; Value = "-4.96%" }
    @{ Cell = "D38"; Value = "0.02672" }
    @{ Cell = "E38"; Value = "11.42%" }
    @{ Cell = "D39"; Value = "0.05604" }
    @{ Cell = "E39"; Value = "7.63%" }
    @{ Cell = "D40"; Value = "0.006300" }
    @{ Cell = "D41"; Value = "0.007691" }
    @{ Cell = "E41"; Value = "-1.43%" }
    @{ Cell = "D42"; Value = "0.1440" }
    @{ Cell = "E42"; Value = "8.63%" }
    @{ Cell = "D43"; Value = "0.007694" }
    @{ Cell = "E43"; Value = "3.79%" }
    @{ Cell = "E44"; Value = "6.74%" }
    @{ Cell = "E45"; Value = "-1.47%" }
    @{ Cell = "D46"; Value = "0.00006990" }
    @{ Cell = "E46"; Value = "6.18%" }
    @{ Cell = "D48"; Value = "0.05351" }
    @{ Cell = "E48"; Value = "18.12%" }
    @{ Cell = "D49"; Value = "0.004000" }
    @{ Cell = "E49"; Value = "-4.82%" }
    @{ Cell = "D50"; Value = "0.00002100" }
    @{ Cell = "D51"; Value = "0.0002000" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
